# ---------------------------------------------------------------------------
# Applies the "Version 2 with station types" edit to the SALB benchmark
# workbook:
#   - num_products on the overview sheet goes from 2 to 1
#   - task_times loses the times_product2 column (now just one product,
#     renamed header "Product_A")
#   - two brand-new sheets are inserted after task_times:
#       station_types (task_ID / Station_A / Station_B)
#       station_costs (Station / Costs)
#   - precedence_relations gains two additional relation rows
#   - the final selection / active tab ends up on precedence_relations
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. overview: num_products 2 -> 1
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("overview")
$wsOverview.Range("B3").Value = 1

# ---------------------------------------------------------------------------
# 2. task_times: drop the times_product2 column (rename happens below, after
#    the station_types headers, to mirror the original authoring order)
# ---------------------------------------------------------------------------
$wsTaskTimes = $wb.Worksheets.Item("task_times")
$wsTaskTimes.Range("C1:C11").ClearContents()

# ---------------------------------------------------------------------------
# 3. Insert the two new sheets right after task_times
# ---------------------------------------------------------------------------
$wsStationTypes = $wb.Worksheets.Add($null, $wsTaskTimes)
$wsStationTypes.Name = "station_types"

$wsStationCosts = $wb.Worksheets.Add($null, $wsStationTypes)
$wsStationCosts.Name = "station_costs"

# --- station_types: task_ID / Station_A / Station_B ------------------------
# (header order below matches the authoring order baked into the workbook's
#  shared-string table: Station_A, Station_B, Product_A, task_ID, ...)
$wsStationTypes.Range("B1").Value = "Station_A"
$wsStationTypes.Range("C1").Value = "Station_B"

$wsTaskTimes.Range("B1").Value = "Product_A"

$wsStationTypes.Range("A1").Value = "task_ID"

$stationTypeData = @(
    @(1, 1, 0),
    @(2, 1, 0),
    @(3, 0, 1),
    @(4, 0, 1),
    @(5, 1, 1),
    @(6, 1, 1),
    @(7, 1, 0),
    @(8, 1, 1),
    @(9, 1, 0),
    @(10, 0, 1)
)
for ($i = 0; $i -lt $stationTypeData.Count; $i++) {
    $r = $i + 2
    $wsStationTypes.Cells.Item($r, 1).Value = $stationTypeData[$i][0]
    $wsStationTypes.Cells.Item($r, 2).Value = $stationTypeData[$i][1]
    $wsStationTypes.Cells.Item($r, 3).Value = $stationTypeData[$i][2]
}

# --- station_costs: Station / Costs -----------------------------------------
$wsStationCosts.Range("A1").Value = "Station"
$wsStationCosts.Range("B1").Value = "Costs"
$wsStationCosts.Range("A2").Value = "Station_A"
$wsStationCosts.Range("B2").Value = 50
$wsStationCosts.Range("A3").Value = "Station_B"
$wsStationCosts.Range("B3").Value = 60

# ---------------------------------------------------------------------------
# 4. precedence_relations: append two new relation rows
# ---------------------------------------------------------------------------
$wsPrecedence = $wb.Worksheets.Item("precedence_relations")
$wsPrecedence.Range("A3").Value = "4;1"
$wsPrecedence.Range("A4").Value = "2;1"

# ---------------------------------------------------------------------------
# 5. Replay the final selection state seen in the authored workbook
# ---------------------------------------------------------------------------
[void]$wsOverview.Select()
[void]$wsOverview.Range("E14").Select()

[void]$wsTaskTimes.Select()
[void]$wsTaskTimes.Range("E18").Select()

[void]$wsStationTypes.Select()
[void]$wsStationTypes.Range("C14").Select()

[void]$wsStationCosts.Select()
[void]$wsStationCosts.Range("B22").Select()

$wsIncompatible = $wb.Worksheets.Item("incompatible_tasks")
[void]$wsIncompatible.Select()
[void]$wsIncompatible.Range("F29").Select()

$wsCompatible = $wb.Worksheets.Item("compatible_tasks")
[void]$wsCompatible.Select()

# precedence_relations ends up the active tab with A13 selected
[void]$wsPrecedence.Select()
[void]$wsPrecedence.Range("A13").Select()
